$wb = $excel.ActiveWorkbook

# --- Sheet "Typography" (sheet1) ---
$wsTypo = $wb.Worksheets.Item("Typography")

# Rows 9 & 10 (Typography_02 / Typography_03): font swisop3.ttf -> JetBrainsMono-Regular.ttf
$wsTypo.Range("C9").Value = "JetBrainsMono-Regular.ttf"
$wsTypo.Range("C10").Value = "JetBrainsMono-Regular.ttf"

# New row 12: Typography_05 using JetBrainsMono-Regular.ttf
$wsTypo.Range("B12").Value = "Typography_05"
$wsTypo.Range("C12").Value = "JetBrainsMono-Regular.ttf"
$wsTypo.Range("D12").Value = 60
$wsTypo.Range("E12").Value = 4
$wsTypo.Range("F12").Value = "?"
$wsTypo.Range("G12").Value = ""
$wsTypo.Range("H12").Value = ""
$wsTypo.Range("I12").Value = ""
$wsTypo.Range("J12").Value = ""

# New row 13: Typography_06 using verdana.ttf
$wsTypo.Range("B13").Value = "Typography_06"
$wsTypo.Range("C13").Value = "verdana.ttf"
$wsTypo.Range("D13").Value = 20
$wsTypo.Range("E13").Value = 4
$wsTypo.Range("F13").Value = "?"
$wsTypo.Range("G13").Value = ""
$wsTypo.Range("H13").Value = ""
$wsTypo.Range("I13").Value = ""
$wsTypo.Range("J13").Value = ""

# --- Sheet "Translation" (sheet2) ---
$wsTrans = $wb.Worksheets.Item("Translation")

# Row 5: typography name now Typography_05
$wsTrans.Range("C5").Value = "Typography_05"

# Row 6: fix double space in battery text
$wsTrans.Range("F6").Value = "BATT: <value> V"

# Row 7: new oil-pressure text (was TCS text)
$wsTrans.Range("F7").Value = "OIL PRESS `n<value>"

# Row 8: example value changed
$wsTrans.Range("F8").Value = "26C"

# Row 12: example value changed
$wsTrans.Range("F12").Value = "12.45"

# Row 13: example value changed
$wsTrans.Range("F13").Value = "80"

# Row 14 now carries what used to be row 15's id/typography, with a new GB example
$wsTrans.Range("B14").Value = "SingleUseId13"
$wsTrans.Range("C14").Value = "Typography_05"
$wsTrans.Range("F14").Value = "12345"

# Row 15 now carries what used to be row 16's content
$wsTrans.Range("B15").Value = "SingleUseId14"
$wsTrans.Range("C15").Value = "Typography_01"
$wsTrans.Range("F15").Value = "2"

# Row 16 is now empty (old row 16 content moved up to row 15)
$wsTrans.Range("B16:F16").ClearContents()
